# Auto-generated edit script: updates Leve profit-calculation columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets, per the scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6495042.5
$ws.Range("I33").Value = 1759.1666
$ws.Range("K33").Value = 1759.1666
$ws.Range("M33").Value = -1530.1666

$ws.Range("H64").Value = 3693.182
$ws.Range("I64").Value = 3430.8
$ws.Range("J64").Value = 4255.4287
$ws.Range("K64").Value = 3430.8
$ws.Range("L64").Value = 4255.4287
$ws.Range("M64").Value = -3182.8
$ws.Range("N64").Value = -4751.4287

$ws.Range("H67").Value = 3693.182
$ws.Range("I67").Value = 3430.8
$ws.Range("J67").Value = 4255.4287
$ws.Range("K67").Value = 3430.8
$ws.Range("L67").Value = 4255.4287
$ws.Range("M67").Value = -2572.8
$ws.Range("N67").Value = -5971.4287

$ws.Range("H113").Value = 5583.5713
$ws.Range("I113").Value = 4523.625
$ws.Range("J113").Value = 6996.8335
$ws.Range("K113").Value = 4523.625
$ws.Range("L113").Value = 6996.8335
$ws.Range("M113").Value = -1269.625
$ws.Range("N113").Value = -13504.8335

$ws.Range("H132").Value = 1603.3684
$ws.Range("I132").Value = 1136.8889
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 3410.6667
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -880.6666999999998
$ws.Range("N132").Value = -35060

$ws.Range("H138").Value = 4110.369
$ws.Range("I138").Value = 1202.4359
$ws.Range("J138").Value = 6630.5776
$ws.Range("K138").Value = 3607.3077
$ws.Range("L138").Value = 19891.7328
$ws.Range("M138").Value = 1532.6923
$ws.Range("N138").Value = -30171.7328

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1543.3334
$ws.Range("I2").Value = 1750.7
$ws.Range("J2").Value = 506.5
$ws.Range("K2").Value = 1750.7
$ws.Range("L2").Value = 506.5
$ws.Range("M2").Value = -1637.7
$ws.Range("N2").Value = -732.5

$ws.Range("H32").Value = 7317.6943
$ws.Range("I32").Value = 7084.92
$ws.Range("J32").Value = 7846.727
$ws.Range("K32").Value = 7084.92
$ws.Range("L32").Value = 7846.727
$ws.Range("M32").Value = -6797.92
$ws.Range("N32").Value = -8420.726999999999

$ws.Range("H45").Value = 14884
$ws.Range("I45").Value = 25937.334
$ws.Range("J45").Value = 1620
$ws.Range("K45").Value = 25937.334
$ws.Range("L45").Value = 1620
$ws.Range("M45").Value = -25560.334
$ws.Range("N45").Value = -2374

$ws.Range("H61").Value = 1962.2449
$ws.Range("I61").Value = 1919.0968
$ws.Range("J61").Value = 2036.5555
$ws.Range("K61").Value = 1919.0968
$ws.Range("L61").Value = 2036.5555
$ws.Range("M61").Value = -1707.0968
$ws.Range("N61").Value = -2460.5555

$ws.Range("H102").Value = 3377805
$ws.Range("I102").Value = 3377805
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3377805
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3376183
$ws.Range("N102").ClearContents()

$ws.Range("H116").Value = 1543.3334
$ws.Range("I116").Value = 1750.7
$ws.Range("J116").Value = 506.5
$ws.Range("K116").Value = 1750.7
$ws.Range("L116").Value = 506.5
$ws.Range("M116").Value = 543.3
$ws.Range("N116").Value = -5094.5

$ws.Range("H122").Value = 2138035.2
$ws.Range("I122").Value = 3663974.8
$ws.Range("J122").Value = 1719.8
$ws.Range("K122").Value = 10991924.4
$ws.Range("L122").Value = 5159.4
$ws.Range("M122").Value = -10989474.4
$ws.Range("N122").Value = -10059.4

$ws.Range("H132").Value = 2165.3823
$ws.Range("I132").Value = 1963.7407
$ws.Range("J132").Value = 2943.1428
$ws.Range("K132").Value = 5891.2221
$ws.Range("L132").Value = 8829.428400000001
$ws.Range("M132").Value = -3361.2221
$ws.Range("N132").Value = -13889.4284

$ws.Range("H136").Value = 1962.2449
$ws.Range("I136").Value = 1919.0968
$ws.Range("J136").Value = 2036.5555
$ws.Range("K136").Value = 5757.2904
$ws.Range("L136").Value = 6109.666499999999
$ws.Range("M136").Value = -3207.2904
$ws.Range("N136").Value = -11209.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1543.3334
$ws.Range("I3").Value = 1750.7
$ws.Range("J3").Value = 506.5
$ws.Range("K3").Value = 1750.7
$ws.Range("L3").Value = 506.5
$ws.Range("M3").Value = -1636.7
$ws.Range("N3").Value = -734.5

$ws.Range("H107").Value = 1709.9474
$ws.Range("I107").Value = 1485
$ws.Range("J107").Value = 2339.8
$ws.Range("K107").Value = 1485
$ws.Range("L107").Value = 2339.8
$ws.Range("M107").Value = 435
$ws.Range("N107").Value = -6179.8

$ws.Range("H132").Value = 32260
$ws.Range("J132").Value = 32260
$ws.Range("L132").Value = 32260
$ws.Range("N132").Value = -42380

$ws.Range("H134").Value = 3527.4666
$ws.Range("I134").Value = 3452
$ws.Range("J134").Value = 3577.7778
$ws.Range("K134").Value = 10356
$ws.Range("L134").Value = 10733.3334
$ws.Range("M134").Value = -7821
$ws.Range("N134").Value = -15803.3334

$ws.Range("H135").Value = 149990
$ws.Range("J135").Value = 149990
$ws.Range("L135").Value = 149990
$ws.Range("N135").Value = -160130

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3444.9033
$ws.Range("I134").Value = 3049.2
$ws.Range("J134").Value = 4164.364
$ws.Range("K134").Value = 9147.599999999999
$ws.Range("L134").Value = 12493.092
$ws.Range("M134").Value = -6612.599999999999
$ws.Range("N134").Value = -17563.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4465.25
$ws.Range("I63").Value = 1305.5
$ws.Range("J63").Value = 7625
$ws.Range("K63").Value = 3916.5
$ws.Range("L63").Value = 22875
$ws.Range("M63").Value = -3167.5
$ws.Range("N63").Value = -24373

$ws.Range("H66").Value = 4465.25
$ws.Range("I66").Value = 1305.5
$ws.Range("J66").Value = 7625
$ws.Range("K66").Value = 11749.5
$ws.Range("L66").Value = 68625
$ws.Range("M66").Value = -8005.5
$ws.Range("N66").Value = -76113

$ws.Range("H68").Value = 3141.6
$ws.Range("I68").Value = 4189.8667
$ws.Range("K68").Value = 12569.6001
$ws.Range("M68").Value = -11758.6001

$ws.Range("H71").Value = 3141.6
$ws.Range("I71").Value = 4189.8667
$ws.Range("K71").Value = 37708.8003
$ws.Range("M71").Value = -33652.8003

$ws.Range("H113").Value = 1463887.1
$ws.Range("I113").Value = 2083782.9
$ws.Range("J113").Value = 588740.25
$ws.Range("K113").Value = 6251348.699999999
$ws.Range("L113").Value = 1766220.75
$ws.Range("M113").Value = -6249178.699999999
$ws.Range("N113").Value = -1770560.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1057.5
$ws.Range("I107").Value = 491.4762
$ws.Range("J107").Value = 2755.5715
$ws.Range("K107").Value = 491.4762
$ws.Range("L107").Value = 2755.5715
$ws.Range("M107").Value = 1428.5238
$ws.Range("N107").Value = -6595.5715

$ws.Range("H113").Value = 83334350
$ws.Range("I113").Value = 250000900
$ws.Range("J113").Value = 1078.25
$ws.Range("K113").Value = 250000900
$ws.Range("L113").Value = 1078.25
$ws.Range("M113").Value = -249998730
$ws.Range("N113").Value = -5418.25

$ws.Range("H132").Value = 2902.434
$ws.Range("I132").Value = 3349.7144
$ws.Range("J132").Value = 2608.9062
$ws.Range("K132").Value = 10049.1432
$ws.Range("L132").Value = 7826.7186
$ws.Range("M132").Value = -7519.143199999999
$ws.Range("N132").Value = -12886.7186

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 100003570
$ws.Range("I40").Value = 200003040
$ws.Range("J40").Value = 4091
$ws.Range("K40").Value = 200003040
$ws.Range("L40").Value = 4091
$ws.Range("M40").Value = -200002904
$ws.Range("N40").Value = -4363
